# SubSequentRenewal45_UT_SS.xlsx update
# - Bump the SYMBOL year tag from SYMBOL_2000 to SYMBOL_2017 on every data row
# - Replace the old placeholder symbol codes (C / N / K) in the
#   BI_SYMBOL/PD_SYMBOL/UM_SYMBOL/MP_SYMBOL columns (AC:AF) with the new
#   per-row codes (BI00x/PD00x/UM00x/MP00x)
# - Make rows 2 & 3 a touch shorter (13.2pt, explicit custom height)
# - Move the active selection to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: SYMBOL_2000 -> SYMBOL_2017 for all 4 data rows ---
$ws.Range("B2").Value = "SYMBOL_2017"
$ws.Range("B3").Value = "SYMBOL_2017"
$ws.Range("B4").Value = "SYMBOL_2017"
$ws.Range("B5").Value = "SYMBOL_2017"

# --- Row 2: AC..AF -> BI001/PD001/UM001/MP001 ---
$ws.Range("AC2").Value = "BI001"
$ws.Range("AD2").Value = "PD001"
$ws.Range("AE2").Value = "UM001"
$ws.Range("AF2").Value = "MP001"

# --- Row 3: AC..AF -> BI002/PD002/UM002/MP002 ---
$ws.Range("AC3").Value = "BI002"
$ws.Range("AD3").Value = "PD002"
$ws.Range("AE3").Value = "UM002"
$ws.Range("AF3").Value = "MP002"

# --- Row 4: AC..AF -> BI003/PD003/UM003/MP003 ---
$ws.Range("AC4").Value = "BI003"
$ws.Range("AD4").Value = "PD003"
$ws.Range("AE4").Value = "UM003"
$ws.Range("AF4").Value = "MP003"

# --- Row 5: AC..AF -> BI004/PD004/UM004/MP004 ---
$ws.Range("AC5").Value = "BI004"
$ws.Range("AD5").Value = "PD004"
$ws.Range("AE5").Value = "UM004"
$ws.Range("AF5").Value = "MP004"

# --- Row heights: rows 2 & 3 become 13.2pt custom height ---
$ws.Rows.Item(2).RowHeight = 13.2
$ws.Rows.Item(3).RowHeight = 13.2

# --- Move the selection to B9 (matches the saved sheetView state) ---
$ws.Range("B9").Select()
